# The bullet "Users should be asked questions 1 after the other." lives in
# the text frame of shape "Rectangle 18" (9th shape in the COM Shapes
# collection for this slide), as the first paragraph.
#
# The edit splits the single run of that paragraph into three runs so the
# word "should" can be given its own run/formatting (Consolas font) while
# keeping the rest of the sentence's original formatting (color 373A36,
# "inherit" latin typeface) intact.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(9)

$paragraph = $shape.TextFrame.TextRange.Paragraphs(1)

# "Users should be asked questions 1 after the other."
#  123456789012345678
# "Users " -> chars 1-6 (length 6)
# "should" -> chars 7-12 (length 6)
# " be asked questions 1 after the other." -> the remainder

$word = $paragraph.Characters(7, 6)
$word.Text = "should"
$word.Font.Name = "Consolas"
